$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(3).Copy() | Out-Null
$ws.Rows.Item(4).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown) | Out-Null
$ws.Rows.Item(3).Copy() | Out-Null
$ws.Rows.Item(4).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown) | Out-Null
$excel.CutCopyMode = 0

Write-Host "done"
